$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 746, pushing old rows 746-793 down to 747-794.
$ws.Rows(746).Insert()

# Populate the newly inserted row 746 with the new record
# (columns A,B,C,E,F,G,H,I,J match every other row in this block; K,L,N,O,P,R
# are copied from the row that used to be 746 and is now 747).
$ws.Range("A746").Value = 10
$ws.Range("B746").Value = "Vega Modelo de Temuco"
$ws.Range("C746").Value = "La Araucanía"
$ws.Range("D746").Value = 44946
$ws.Range("E746").Value = 9
$ws.Range("F746").Value = "Fruta"
$ws.Range("G746").Value = 100102
$ws.Range("H746").Value = "Cítricos"
$ws.Range("I746").Value = 100102004
$ws.Range("J746").Value = "Mandarina"
$ws.Range("K746").Value = "Murcott"
$ws.Range("L746").Value = "Primera"
$ws.Range("M746").Value = 118
$ws.Range("N746").Value = 12000
$ws.Range("O746").Value = 12000
$ws.Range("P746").Value = 12000
$ws.Range("Q746").Value = "$/bandeja 10 kilos"
$ws.Range("R746").Value = "Región de O'Higgins"
$ws.Range("S746").Value = 1200
$ws.Range("T746").Value = 10

# Match the date style used by the other rows in column D (style index 2
# in the original file: a yyyy-mm-dd hh:mm:ss number format).
$ws.Range("D746").NumberFormat = $ws.Range("D747").NumberFormat
